$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.298.32"
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("E2").Value = "  +0.04%  "

$ws.Range("D3").Value = "'3.679.17"
$ws.Range("D3").Style = $ws.Range("C3").Style
$ws.Range("E3").Value = "  -0.05%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'675.97"
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = "  -1.22%  "

$ws.Range("D6").Value = "'158.39"
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = "  -2.45%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.23%  "

$ws.Range("E9").Value = "  -1.40%  "

$ws.Range("D10").Value = "'6.93"
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = "  -6.20%  "

$ws.Range("E11").Value = "  -2.07%  "

$ws.Range("D12").Value = "'0.0000232"
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = "  -3.41%  "

$ws.Range("D13").Value = "'4.300.13"
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = "  -0.04%  "

$ws.Range("D14").Value = "'32.33"
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = "  -3.82%  "

$ws.Range("D15").Value = "'3.670.86"
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("D16").Value = "'69.268.19"
$ws.Range("D16").Style = $ws.Range("C16").Style
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("D17").Value = "'0.114"
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = "  +1.53%  "

$ws.Range("D18").Value = "'16.03"
$ws.Range("D18").Style = $ws.Range("C18").Style
$ws.Range("E18").Value = "  -2.04%  "

$ws.Range("E19").Value = "  -3.09%  "

$ws.Range("D20").Value = "'468.73"
$ws.Range("D20").Style = $ws.Range("C20").Style
$ws.Range("E20").Value = "  -2.93%  "

$ws.Range("D21").Value = "'9.95"
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = "  +0.09%  "

$ws.Range("E22").Value = "  -2.53%  "

$ws.Range("D23").Value = "'79.77"
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").Value = "'3.823.44"
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  -5.77%  "

$ws.Range("D27").Value = "'10.89"
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = "  -4.53%  "

$ws.Range("D28").Value = "'9.06"
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = "  -4.67%  "

$ws.Range("E29").Value = "  -1.48%  "

$ws.Range("E30").Value = "  -4.42%  "

$ws.Range("D31").Value = "'6.61"
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = "  -3.37%  "

$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").Value = "'26.93"
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = "  -0.78%  "

$ws.Range("E34").Value = "  -4.69%  "

$ws.Range("D35").Value = "'3.671.03"
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").Value = "'0.158"
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = "  -5.33%  "

$ws.Range("D37").Value = "'8.21"
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = "  -3.13%  "

$ws.Range("E38").Value = "  -1.79%  "

$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "'2.21"
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = "  -5.28%  "

$ws.Range("D42").Value = "'0.0904"
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").Value = "'172.05"
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = "  +5.39%  "

$ws.Range("D44").Value = "'0.941"
$ws.Range("D44").Style = $ws.Range("C44").Style
$ws.Range("E44").Value = "  -0.87%  "

$ws.Range("D45").Value = "'47.66"
$ws.Range("D45").Style = $ws.Range("C45").Style
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("D46").Value = "'0.000277"
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = "  -3.20%  "

$ws.Range("E47").Value = "  -5.28%  "

$ws.Range("D48").Value = "'27.59"
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = "  -7.01%  "

$ws.Range("D49").Value = "'1.28"
$ws.Range("D49").Style = $ws.Range("C49").Style
$ws.Range("E49").Value = "  -4.16%  "

$ws.Range("E50").Value = "  -2.79%  "

$ws.Range("E51").Value = "  -2.87%  "
